# Combine Florenceville and Bristol prior to 2009 (#4)
#
# Starting with 2009, Florenceville and Bristol were amalgamated into
# Florenceville-Bristol. For data consistency purposes, the standalone
# "Bristol" and "Florenceville" rows are removed from this historical
# sheet (the data has already been manually combined into the existing
# "Florenceville-Bristol" row elsewhere in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet holds an Excel Table ("Frame0") with Municipality in column A
# and Policing Provider in column B. Find and delete the rows whose
# Municipality is exactly "Bristol" or "Florenceville" (but not
# "Florenceville-Bristol"). Deleting entire rows from the bottom up keeps
# earlier row indices valid and causes the table/used range to shrink
# automatically.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

$rowsToDelete = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($name -eq "Bristol" -or $name -eq "Florenceville") {
        $rowsToDelete += $r
    }
}

# Delete from bottom-most row first so indices of rows still to be
# deleted remain valid.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
